$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "57.757.45"
$ws.Range("E2").Value = "  -5.65%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.909.02"
$ws.Range("E3").Value = "  -3.48%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
Set-TextValue "D5" "550.17"
$ws.Range("E5").Value = "  -3.24%  "

# Row 6 - Solana
Set-TextValue "D6" "123.42"
$ws.Range("E6").Value = "  -4.25%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.06%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "2.905.98"
$ws.Range("E8").Value = "  -3.50%  "

# Row 9 - XRP
Set-TextValue "D9" "0.496"
$ws.Range("E9").Value = "  -0.14%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.125"
$ws.Range("E10").Value = "  -6.65%  "

# Row 11 - Toncoin
Set-TextValue "D11" "4.77"
$ws.Range("E11").Value = "  -7.88%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.440"
$ws.Range("E12").Value = "  +2.53%  "

# Row 13 - ShibaInu
Set-TextValue "D13" "0.0000213"
$ws.Range("E13").Value = "  -4.48%  "

# Row 14 - Avalanche
Set-TextValue "D14" "32.46"
$ws.Range("E14").Value = "  -1.01%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  +1.49%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "3.388.53"
$ws.Range("E16").Value = "  -3.50%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.908.00"
$ws.Range("E17").Value = "  -3.53%  "

# Row 18 - Polkadot
$ws.Range("E18").Value = "  +6.50%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "57.774.27"
$ws.Range("E19").Value = "  -5.77%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "410.37"
$ws.Range("E20").Value = "  -6.32%  "

# Row 21 - Chainlink
$ws.Range("E21").Value = "  -1.63%  "

# Row 22 - Polygon
Set-TextValue "D22" "0.677"
$ws.Range("E22").Value = "  +2.52%  "

# Row 23 - Uniswap
Set-TextValue "D23" "6.87"
$ws.Range("E23").Value = "  -3.64%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue "D24" "12.90"
$ws.Range("E24").Value = "  +3.04%  "

# Row 25 - Litecoin
Set-TextValue "D25" "77.32"
$ws.Range("E25").Value = "  -2.02%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.11%  "

# Row 27 - FirstDigitalUSD
$ws.Range("E27").Value = "  -0.11%  "

# Row 28 - PancakeSwap
$ws.Range("E28").Value = "  -1.54%  "

# Row 29 - RenderToken->ImmutableX
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D29" "1.95"
$ws.Range("E29").Value = "  +3.84%  "

# Row 30 - ImmutableX->RenderToken
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D30" "7.28"
$ws.Range("E30").Value = "  +1.00%  "

# Row 31 - NEARProtocol
Set-TextValue "D31" "6.09"
$ws.Range("E31").Value = "  -1.83%  "

# Row 32 - EthereumClassic
Set-TextValue "D32" "24.74"
$ws.Range("E32").Value = "  -3.03%  "

# Row 33 - Hedera
Set-TextValue "D33" "0.0984"
$ws.Range("E33").Value = "  +4.51%  "

# Row 34 - Mantle
Set-TextValue "D34" "0.915"
$ws.Range("E34").Value = "  -4.20%  "

# Row 35 - Filecoin
Set-TextValue "D35" "5.40"
$ws.Range("E35").Value = "  -2.89%  "

# Row 36 - Stacks
Set-TextValue "D36" "2.01"
$ws.Range("E36").Value = "  -11.23%  "

# Row 37 - OKB
Set-TextValue "D37" "48.18"
$ws.Range("E37").Value = "  -3.80%  "

# Row 38 - Cosmos
$ws.Range("E38").Value = "  +9.46%  "

# Row 39 - PEPE
$ws.Range("D39").Value = "0.0₃0624"
$ws.Range("E39").Value = "  -8.35%  "

# Row 40 - VeChain
Set-TextValue "D40" "0.0346"

# Row 41 - Kaspa
Set-TextValue "D41" "0.106"
$ws.Range("E41").Value = "  -1.30%  "

# Row 42 - Maker
$ws.Range("D42").Value = "2.629.35"
$ws.Range("E42").Value = "  -0.66%  "

# Row 43 - Bittensor
Set-TextValue "D43" "362.07"
$ws.Range("E43").Value = "  -2.91%  "

# Row 44 - dogwifhat
Set-TextValue "D44" "2.40"
$ws.Range("E44").Value = "  -1.25%  "

# Row 45 - USDe
Set-TextValue "D45" "0.998"
$ws.Range("E45").Value = "  -0.03%  "

# Row 46 - Monero
Set-TextValue "D46" "120.01"
$ws.Range("E46").Value = "  +0.62%  "

# Row 47 - TheGraph
Set-TextValue "D47" "0.229"
$ws.Range("E47").Value = "  -2.69%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  +1.12%  "

# Row 49 - Fetch.AI
Set-TextValue "D49" "1.95"
$ws.Range("E49").Value = "  -0.84%  "

# Row 50 - InjectiveProtocol
$ws.Range("E50").Value = "  -3.09%  "

# Row 51 - ThetaToken
$ws.Range("E51").Value = "  -3.03%  "
